# Tennessee overview workbook — convert numeric "count" cells to literal text
# (to match values pasted in from a separately formatted export), add a
# missing "Total" row to the County sheet, and fix the placeholder zero row
# for Fentress County (which had no filers) to show percent/currency text.
#
# NOTE: this engine's PowerShell subset only binds *positional* function
# parameters, not -Named ones, so helpers below take plain positional args.

function Set-TextValue {
    param($Cell, $Text)
    # Force the cell into Text format so a numeric-looking string (e.g. "24"
    # or "2,030") is stored verbatim instead of being re-parsed as a number,
    # then drop back to the workbook's default style so no stray number
    # format is left behind on the cell.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (No. of 990 Filers w/ Gov Grants) 2030 -> "2,030"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Cells.Item(2, 1) "2,030"

# ---------------------------------------------------------------------
# Sheet "County": convert B2:B90 counts to text, fix the Fentress County
# placeholder row (91), and append a new Total row (92).
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

for ($r = 2; $r -le 90; $r++) {
    $cell = $wsCounty.Cells.Item($r, 2)
    $val = $cell.Value2
    $valText = [string]$val
    Set-TextValue $cell $valText
}

# Row 91 (Fentress County) — was all literal "0"s, now shows formatted
# percent/currency placeholders.
Set-TextValue $wsCounty.Cells.Item(91, 2) "0.00%"
Set-TextValue $wsCounty.Cells.Item(91, 3) "`$0"
Set-TextValue $wsCounty.Cells.Item(91, 4) "0.00%"
Set-TextValue $wsCounty.Cells.Item(91, 5) "0.00%"
Set-TextValue $wsCounty.Cells.Item(91, 6) "0.00%"

# New row 92 — statewide Total, matching the Overall sheet figures.
Set-TextValue $wsCounty.Cells.Item(92, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(92, 2) "2,030"
Set-TextValue $wsCounty.Cells.Item(92, 3) "`$3,759,003,111"
Set-TextValue $wsCounty.Cells.Item(92, 4) "8.54%"
Set-TextValue $wsCounty.Cells.Item(92, 5) "-16.99%"
Set-TextValue $wsCounty.Cells.Item(92, 6) "69.16%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": convert B2:B10 counts + B11 Total to text.
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
for ($r = 2; $r -le 10; $r++) {
    $cell = $wsCd.Cells.Item($r, 2)
    $val = $cell.Value2
    $valText = [string]$val
    Set-TextValue $cell $valText
}
Set-TextValue $wsCd.Cells.Item(11, 2) "2,030"

# ---------------------------------------------------------------------
# Sheet "Size": convert B2:B7 counts + B8 Total to text.
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
for ($r = 2; $r -le 7; $r++) {
    $cell = $wsSize.Cells.Item($r, 2)
    $val = $cell.Value2
    $valText = [string]$val
    Set-TextValue $cell $valText
}
Set-TextValue $wsSize.Cells.Item(8, 2) "2,030"

# ---------------------------------------------------------------------
# Sheet "Subsector": convert B2:B13 counts + B14 Total to text.
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
for ($r = 2; $r -le 13; $r++) {
    $cell = $wsSub.Cells.Item($r, 2)
    $val = $cell.Value2
    $valText = [string]$val
    Set-TextValue $cell $valText
}
Set-TextValue $wsSub.Cells.Item(14, 2) "2,030"

Write-Host "Tennessee overview text-format edits applied."
